# FlixelRL-198: 武器パラメータの調整 (weapon/armor name + atk/def rebalance)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("item_equipment")

# --- Weapon names (column C, rows 3-9 = WEAPON1..WEAPON7) ---
$ws.Range("C3").Value = "竹刀"
$ws.Range("C4").Value = "木刀"
$ws.Range("C5").Value = "ダガー"
$ws.Range("C6").Value = "レイピア"
$ws.Range("C7").Value = "三日月刀"
$ws.Range("C8").Value = "妖刀ムラマサ"
$ws.Range("C9").Value = "ライトセーバー"

# --- Weapon atk (column D, rows 3-9) ---
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 2
$ws.Range("D5").Value = 3
$ws.Range("D6").Value = 5
$ws.Range("D7").Value = 10
$ws.Range("D8").Value = 20
$ws.Range("D9").Value = 25

# --- Armor names (column C, rows 23-29 = ARMOR1..ARMOR7) ---
$ws.Range("C23").Value = "ローブ"
$ws.Range("C24").Value = "毛皮の鎧"
$ws.Range("C25").Value = "鎖かたびら"
$ws.Range("C26").Value = "エルフの鎧"
$ws.Range("C27").Value = "鋼鉄の鎧"
$ws.Range("C28").Value = "シルバーアーマー"
$ws.Range("C29").Value = "プラチナメイル"

# --- Armor def (column E, rows 23-29) ---
$ws.Range("E23").Value = 1
$ws.Range("E24").Value = 2
$ws.Range("E25").Value = 3
$ws.Range("E26").Value = 5
$ws.Range("E27").Value = 10
$ws.Range("E28").Value = 20
$ws.Range("E29").Value = 25
